$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting existing rows 3-5 down to 4-6.
$ws.Rows.Item(3).Insert()

# The newly inserted row 3 is currently blank; populate it with the
# weekly data row (same constant columns as neighboring rows, new
# date/volume/price values).
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Macroferia Regional de Talca"
$ws.Range("C3").Value = "Maule"
$ws.Range("D3").Value = 45274
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101004
$ws.Range("J3").Value = "Frambuesa"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 3500
$ws.Range("T3").Value = 2
